# Slide 1, Shape 1 ("Subtitle 2") holds the bulleted "Slides Content" list.
# The 6th paragraph's wording changes from "...can be incorporated into a
# report." to "...can be integrated into a report."
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tf = $shape.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(6)

# First overwrite with unrelated placeholder text, then set the final
# wording. This avoids the host's internal char-level diff matching common
# substrings between "incorporated" and "integrated" (e.g. "in"/"rated"),
# which would otherwise split the single <a:r> run into several runs. The
# two-step assignment guarantees the paragraph keeps exactly one run, with
# its original rPr (font/formatting) untouched, matching the target edit.
$para.Text = "XQZJKVWPLM"
$para.Text = "Incorporate a paragraph to show how these results can be integrated into a report."
